$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AI2").Value = 1970.676132921329
$ws.Range("AI3").Value = 5451.591659061644
$ws.Range("AI4").Value = 8398.79062580929
$ws.Range("AI5").Value = 15120.55495915894
$ws.Range("AI6").Value = 21592.12108922844
$ws.Range("AI7").Value = 27304.98484780873
$ws.Range("AI8").Value = 32360.63244316275
$ws.Range("AI9").Value = 36839.96172348834
$ws.Range("AI10").Value = 40807.87515397542
$ws.Range("AI11").Value = 44316.61197830974
$ws.Range("AI12").Value = 47408.19017166216
$ws.Range("AI13").Value = 50116.20432672538
$ws.Range("AI14").Value = 52467.14487454593
$ws.Range("AI15").Value = 54481.35000948822
$ws.Range("AI16").Value = 56173.6641139854
$ws.Range("AI17").Value = 57553.84892010906
$ws.Range("AI18").Value = 58626.77178063371
$ws.Range("AI19").Value = 59392.37602664282
$ws.Range("AI20").Value = 59845.41860986569
$ws.Range("AI21").Value = 59974.93694152016
$ws.Range("AI22").Value = 59763.37586622223
$ws.Range("AI23").Value = 31040.20166374412
$ws.Range("AI24").Value = 53355.25426620171
$ws.Range("AI25").Value = 67308.01980792976
$ws.Range("AI26").Value = 68100.93149641833
$ws.Range("AI27").Value = 5932.181724629349
$ws.Range("AI28").Value = 42411.46612653913
$ws.Range("AI29").Value = 67399.13600942053
$ws.Range("AI30").Value = 68056.72395030156
$ws.Range("AI31").Value = 8430.712381303023
$ws.Range("AI32").Value = 35061.06583760159
$ws.Range("AI33").Value = 5045.932645310533
$ws.Range("AI34").Value = 28804.4628332987
